$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13, pushing existing rows 13-32 down to 14-33.
$ws.Rows.Item(13).Insert()

# Populate the newly inserted row 13 with the new record.
$ws.Cells.Item(13, 1).Value = 4
$ws.Cells.Item(13, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(13, 3).Value = "Los Lagos"
$ws.Cells.Item(13, 4).Value = 45177
$ws.Cells.Item(13, 4).NumberFormat = $ws.Cells.Item(14, 4).NumberFormat
$ws.Cells.Item(13, 5).Value = 10
$ws.Cells.Item(13, 6).Value = 100112035
$ws.Cells.Item(13, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(13, 8).Value = "Sin especificar"
$ws.Cells.Item(13, 9).Value = "Primera"
$ws.Cells.Item(13, 10).Value = 120
$ws.Cells.Item(13, 11).Value = 26000
$ws.Cells.Item(13, 12).Value = 26000
$ws.Cells.Item(13, 13).Value = 26000
$ws.Cells.Item(13, 14).Value = "`$/malla 15 kilos"
$ws.Cells.Item(13, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(13, 16).Value = 1733
$ws.Cells.Item(13, 17).Value = 15
$ws.Cells.Item(13, 18).Value = "Hortaliza"
